$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("readData")

# New header + data cells, written column by column so the shared-string
# table is built in the same interleaved order as the target workbook.
$ws.Range("D1").Value = "Extent report name"
$ws.Range("D2").Value = "Make My Trip Hackathon "

$ws.Range("E1").Value = "Login Id"
$ws.Range("E2").Value = "bughunterss01@gmail.com"

$ws.Range("F1").Value = "Login Password"
$ws.Range("F2").Value = "Bughunter$6"

$ws.Range("G1").Value = "PopUp Title"
$ws.Range("G2").Value = "Login/Signup for Best Prices"

$ws.Range("H1").Value = "Departure city input"
$ws.Range("H2").Value = "Delhi"

$ws.Range("I1").Value = "Arrival City Input"
$ws.Range("I2").Value = "Manali"

# Bold the header row (A1:I1) and keep wrap text
$headerRange = $ws.Range("A1:I1")
$headerRange.Font.Bold = $true
$headerRange.WrapText = $true

# Adjust header row height
$ws.Rows.Item(1).RowHeight = 43.2

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection
[void]$ws.Range("I10").Select()

Write-Host "done"
